$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Step 2 row (row 3) - Log in text and expected result
$ws.Range("C3").Value = "Step 2: Log in as a user with the appropriate role"
$ws.Range("D3").Value = "I am redirected to the user's dashboard"

# Update Step 3 row (row 4) - Navigate to Team KPI page
$ws.Range("C4").Value = "Step 3: Go to the ""Team KPI"" page "
$ws.Range("D4").Value = "I am redirected to a page of KPIs of people on my team"

# Update Step 4 row (row 5) - Delete a KPI belonging to someone on my team
$ws.Range("C5").Value = "Step 4:Delete one that belongs to someone on my team"
$ws.Range("D5").Value = "The data is removed from the database."

# Add Step 5 row (row 6) - Try to delete a kpi about me
$ws.Range("C6").Value = "Step 5: While logged in try to delete a kpi about me"
$ws.Range("D6").Value = "I am denied access to this"

# Add Step 6 row (row 7) - Try to delete a kpi of someone not under my team
$ws.Range("C7").Value = "Step 6: Try to delete a kpi of someone who is not under my team"
$ws.Range("D7").Value = "I am denied access to this"

# Update selection to D4
$ws.Range("D4").Select()
